$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_VS")
$ws.Name = "CRF_VS"
